$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 9 (old) no longer exists; row 8 loses its C/D entries ---
$ws.Range("C8:D8").ClearContents()
$ws.Range("A9:F9").ClearContents()

# --- New semester block: Fall 2023 / Spring 2023 / Summer 2023 (rows 12-14, total row 20) ---
$ws.Range("A12").Value = "Fall 2023"
$ws.Range("B12").Value = "Credits"
$ws.Range("C12").Value = "Spring 2023"
$ws.Range("D12").Value = "Credits"
$ws.Range("E12").Value = "Summer 2023"
$ws.Range("F12").Value = "Credits"

$ws.Range("A13").Value = "CPSC 4176"
$ws.Range("B13").Value = 3

$ws.Range("A14").Value = "CPSC 4000"
$ws.Range("B14").Value = 0

$ws.Range("A20").Value = "Total"
$ws.Range("B20").Formula = "=SUM(B13:B19)"
$ws.Range("C20").Value = "Total"
$ws.Range("D20").Formula = "=SUM(D13:D19)"
$ws.Range("E20").Value = "Total"
$ws.Range("F20").Formula = "=SUM(F13:F19)"

# --- New semester block: Fall 2024 / Spring 2024 / Summer 2024 (header row 21, total row 29) ---
$ws.Range("A21").Value = "Fall 2024"
$ws.Range("B21").Value = "Credits"
$ws.Range("C21").Value = "Spring 2024"
$ws.Range("D21").Value = "Credits"
$ws.Range("E21").Value = "Summer 2024"
$ws.Range("F21").Value = "Credits"

$ws.Range("A29").Value = "Total"
$ws.Range("B29").Formula = "=SUM(B22:B28)"
$ws.Range("C29").Value = "Total"
$ws.Range("D29").Formula = "=SUM(D22:D28)"
$ws.Range("E29").Value = "Total"
$ws.Range("F29").Formula = "=SUM(F22:F28)"

# --- New semester block: Fall 2025 / Spring 2025 / Summer 2025 (header row 30, total row 38) ---
$ws.Range("A30").Value = "Fall 2025"
$ws.Range("B30").Value = "Credits"
$ws.Range("C30").Value = "Spring 2025"
$ws.Range("D30").Value = "Credits"
$ws.Range("E30").Value = "Summer 2025"
$ws.Range("F30").Value = "Credits"

$ws.Range("A38").Value = "Total"
$ws.Range("B38").Formula = "=SUM(B31:B37)"
$ws.Range("C38").Value = "Total"
$ws.Range("D38").Formula = "=SUM(D31:D37)"
$ws.Range("E38").Value = "Total"
$ws.Range("F38").Formula = "=SUM(F31:F37)"
